# Regenerate the s_vals data (filter save games) - update computed columns
# B (TB), C (d2S), D (K), E (IP), G (sum) for each data row. Column A (date)
# and F (Win) remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 },
    @{ Row = 3;  B = 0.1190320826869504; C = 0.306821227259698;   D = 0.1494219747398047; E = 0.4942365360607697; G = 1.069511820747223 },
    @{ Row = 4;  B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 },
    @{ Row = 5;  B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 },
    @{ Row = 6;  B = 1.455362044514542;  C = 10.34677158129881;  D = 0.1494219747398047; E = 10.19245300693656;  G = 22.14400860748972 },
    @{ Row = 7;  B = 0.6606524410359556; C = 0.04071648406533734; D = 0.7527432677738641; E = 0.4942365360607697; G = 1.948348728935927 },
    @{ Row = 8;  B = 1.455362044514542;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 4.358119930609447 },
    @{ Row = 9;  B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 },
    @{ Row = 10; B = 1.455362044514542;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 4.358119930609447 },
    @{ Row = 11; B = 0.2917716402565462; C = 0.306821227259698;   D = 22.3905356188092;   E = 10.19245300693656;  G = 33.181581493262 },
    @{ Row = 12; B = 0.6606524410359556; C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 2.960089034096801 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 7).Value = $item.G
}
